$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed symbol list.
# The new values are numeric-looking strings (plain decimals or "N.NN%"),
# so force text formatting first (NumberFormat = "@") to keep them as text
# instead of letting Excel auto-convert them to numbers / percentages.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "310.08"
$ws.Range("E2").Value = "4.27%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "44.19"
$ws.Range("E3").Value = "6.33%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.088"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07998"
$ws.Range("E5").Value = "5.86%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.460"
$ws.Range("E6").Value = "1.92%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.628"
$ws.Range("E7").Value = "1.64%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "1.070"
$ws.Range("E8").Value = "15.21%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1288"
$ws.Range("E9").Value = "7.02%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1891"
$ws.Range("E10").Value = "2.75%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09186"
$ws.Range("E11").Value = "2.94%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04190"
$ws.Range("E12").Value = "3.78%"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.85%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001305"
$ws.Range("E14").Value = "1.89%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005691"
$ws.Range("E15").Value = "-2.14%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "3.371"
$ws.Range("E17").Value = "0.96%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "2.402"
$ws.Range("E18").Value = "-0.10%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3356"
$ws.Range("E19").Value = "1.33%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "7.980"
$ws.Range("E20").Value = "0.15%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").Value = "-3.33%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3116"
$ws.Range("E22").Value = "3.92%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04172"
$ws.Range("E23").Value = "2.92%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001264"
$ws.Range("E24").Value = "-0.01%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004311"
$ws.Range("E25").Value = "3.83%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001332"
$ws.Range("E26").Value = "8.34%"

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02648"
$ws.Range("E38").Value = "9.57%"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05383"
$ws.Range("E39").Value = "3.17%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005600"
$ws.Range("E40").Value = "-13.96%"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.00%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1406"
$ws.Range("E42").Value = "5.57%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007260"
$ws.Range("E43").Value = "-3.85%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008383"
$ws.Range("E44").Value = "6.97%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3086"
$ws.Range("E45").Value = "-4.14%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006716"
$ws.Range("E46").Value = "-0.98%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-1.27%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05078"
$ws.Range("E48").Value = "9.75%"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003948"
$ws.Range("E49").Value = "-5.97%"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002073"
$ws.Range("E50").Value = "-1.27%"

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001974"
$ws.Range("E51").Value = "-1.27%"
